$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4150
$ws.Range("I76").Value = 4420
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 4420
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -4105
$ws.Range("N76").Value = -3430

$ws.Range("H79").Value = 4150
$ws.Range("I79").Value = 4420
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 4420
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -3328
$ws.Range("N79").Value = -4984

$ws.Range("H132").Value = 8341468
$ws.Range("I132").Value = 16674538
$ws.Range("J132").Value = 8398.549999999999
$ws.Range("K132").Value = 50023614
$ws.Range("L132").Value = 25195.65
$ws.Range("M132").Value = -50021084
$ws.Range("N132").Value = -30255.65

$ws.Range("H137").Value = 1668.9032
$ws.Range("I137").Value = 1458.6666
$ws.Range("J137").Value = 1801.6842
$ws.Range("K137").Value = 4375.9998
$ws.Range("L137").Value = 5405.0526
$ws.Range("M137").Value = -1825.9998
$ws.Range("N137").Value = -10505.0526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5378.3335
$ws.Range("I32").Value = 5067.234
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 5067.234
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -4780.234
$ws.Range("N32").Value = -20574

$ws.Range("H122").Value = 932.53845
$ws.Range("I122").Value = 944.1818
$ws.Range("J122").Value = 868.5
$ws.Range("K122").Value = 2832.5454
$ws.Range("L122").Value = 2605.5
$ws.Range("M122").Value = -382.5454
$ws.Range("N122").Value = -7505.5

$ws.Range("H132").Value = 3969.5833
$ws.Range("I132").Value = 3517.7144
$ws.Range("J132").Value = 4602.2
$ws.Range("K132").Value = 10553.1432
$ws.Range("L132").Value = 13806.6
$ws.Range("M132").Value = -8023.143199999999
$ws.Range("N132").Value = -18866.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 48999.5
$ws.Range("J110").Value = 48999.5
$ws.Range("L110").Value = 48999.5
$ws.Range("N110").Value = -57179.5

$ws.Range("H134").Value = 9732.538
$ws.Range("I134").Value = 1280.4445
$ws.Range("K134").Value = 3841.3335
$ws.Range("M134").Value = -1306.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1112.8889
$ws.Range("I31").Value = 1061.079
$ws.Range("J31").Value = 1394.1428
$ws.Range("K31").Value = 1061.079
$ws.Range("L31").Value = 1394.1428
$ws.Range("M31").Value = -766.079
$ws.Range("N31").Value = -1984.1428

$ws.Range("H34").Value = 1112.8889
$ws.Range("I34").Value = 1061.079
$ws.Range("J34").Value = 1394.1428
$ws.Range("K34").Value = 1061.079
$ws.Range("L34").Value = 1394.1428
$ws.Range("M34").Value = -859.079
$ws.Range("N34").Value = -1798.1428

$ws.Range("H74").Value = 33000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 33000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -34748

$ws.Range("H77").Value = 33000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 99000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -107736

$ws.Range("H94").Value = 1964.3334
$ws.Range("I94").Value = 1943.5
$ws.Range("K94").Value = 1943.5
$ws.Range("M94").Value = -1492.5

$ws.Range("H114").Value = 24866.334
$ws.Range("J114").Value = 24866.334
$ws.Range("L114").Value = 24866.334
$ws.Range("N114").Value = -33544.334

$ws.Range("H132").Value = 9979.615
$ws.Range("I132").Value = 18203.834
$ws.Range("J132").Value = 2930.2856
$ws.Range("K132").Value = 54611.50199999999
$ws.Range("L132").Value = 8790.856800000001
$ws.Range("M132").Value = -52081.50199999999
$ws.Range("N132").Value = -13850.8568

$ws.Range("H134").Value = 31253188
$ws.Range("I134").Value = 4000.1
$ws.Range("K134").Value = 12000.3
$ws.Range("M134").Value = -9465.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6157.6665
$ws.Range("I56").Value = 6157.6665
$ws.Range("K56").Value = 6157.6665
$ws.Range("M56").Value = -5627.6665

$ws.Range("H113").Value = 662.08
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 662.5833
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 1987.7499
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -6327.7499

$ws.Range("H129").Value = 46297204
$ws.Range("I129").Value = 111111730
$ws.Range("J129").Value = 13889946
$ws.Range("K129").Value = 333335190
$ws.Range("L129").Value = 41669838
$ws.Range("M129").Value = -333330190
$ws.Range("N129").Value = -41679838

$ws.Range("H131").Value = 18871214
$ws.Range("J131").Value = 3819.6445
$ws.Range("L131").Value = 11458.9335
$ws.Range("N131").Value = -21538.9335

$ws.Range("H132").Value = 962.25
$ws.Range("I132").Value = 966.3333
$ws.Range("J132").Value = 950
$ws.Range("K132").Value = 8696.9997
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -6166.9997
$ws.Range("N132").Value = -13610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1038
$ws.Range("I22").Value = 974.25
$ws.Range("J22").Value = 1101.75
$ws.Range("K22").Value = 974.25
$ws.Range("L22").Value = 1101.75
$ws.Range("M22").Value = -679.25
$ws.Range("N22").Value = -1691.75

$ws.Range("H27").Value = 1038
$ws.Range("I27").Value = 974.25
$ws.Range("J27").Value = 1101.75
$ws.Range("K27").Value = 974.25
$ws.Range("L27").Value = 1101.75
$ws.Range("M27").Value = -867.25
$ws.Range("N27").Value = -1315.75

$ws.Range("H44").Value = 11900
$ws.Range("J44").Value = 11900
$ws.Range("L44").Value = 11900
$ws.Range("N44").Value = -12812

$ws.Range("H55").Value = 1231.25
$ws.Range("I55").Value = 1175
$ws.Range("J55").Value = 1400
$ws.Range("K55").Value = 1175
$ws.Range("L55").Value = 1400
$ws.Range("M55").Value = -1002
$ws.Range("N55").Value = -1746

$ws.Range("H132").Value = 141088.88
$ws.Range("I132").Value = 55001.5
$ws.Range("J132").Value = 169784.67
$ws.Range("K132").Value = 165004.5
$ws.Range("L132").Value = 509354.01
$ws.Range("M132").Value = -162474.5
$ws.Range("N132").Value = -514414.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
